$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 43863.086
$ws.Range("J2").Value = 167546.5
$ws.Range("L2").Value = 167546.5
$ws.Range("N2").Value = -167772.5
$ws.Range("H11").Value = 143.54546
$ws.Range("I11").Value = 143.54546
$ws.Range("K11").Value = 143.54546
$ws.Range("M11").Value = -3.545459999999991
$ws.Range("H13").Value = 4324.5
$ws.Range("J13").Value = 4432.6665
$ws.Range("L13").Value = 4432.6665
$ws.Range("N13").Value = -4770.6665
$ws.Range("H32").Value = 6712.316
$ws.Range("I32").Value = 5359.4
$ws.Range("K32").Value = 5359.4
$ws.Range("M32").Value = -5033.4
$ws.Range("H33").Value = 1595.92
$ws.Range("I33").Value = 1788.0454
$ws.Range("K33").Value = 1788.0454
$ws.Range("M33").Value = -1559.0454
$ws.Range("H38").Value = 3610
$ws.Range("I38").Value = 47.46154
$ws.Range("J38").Value = 11328.833
$ws.Range("K38").Value = 142.38462
$ws.Range("L38").Value = 33986.499
$ws.Range("M38").Value = 229.61538
$ws.Range("N38").Value = -34730.499
$ws.Range("H40").Value = 4166.5835
$ws.Range("I40").Value = 2555.5557
$ws.Range("J40").Value = 8999.666999999999
$ws.Range("K40").Value = 2555.5557
$ws.Range("L40").Value = 8999.666999999999
$ws.Range("M40").Value = -2380.5557
$ws.Range("N40").Value = -9349.666999999999
$ws.Range("H47").Value = 16877.666
$ws.Range("I47").Value = 133
$ws.Range("J47").Value = 25250
$ws.Range("K47").Value = 133
$ws.Range("L47").Value = 25250
$ws.Range("M47").Value = 839
$ws.Range("N47").Value = -27194
$ws.Range("H62").Value = 4776.2
$ws.Range("I62").Value = 3730.3333
$ws.Range("J62").Value = 6345
$ws.Range("K62").Value = 3730.3333
$ws.Range("L62").Value = 6345
$ws.Range("M62").Value = -3106.3333
$ws.Range("N62").Value = -7593
$ws.Range("H65").Value = 4776.2
$ws.Range("I65").Value = 3730.3333
$ws.Range("J65").Value = 6345
$ws.Range("K65").Value = 18651.6665
$ws.Range("L65").Value = 31725
$ws.Range("M65").Value = -15531.6665
$ws.Range("N65").Value = -37965
$ws.Range("H76").Value = 10221
$ws.Range("I76").Value = 10623.625
$ws.Range("J76").Value = 7000
$ws.Range("K76").Value = 10623.625
$ws.Range("L76").Value = 7000
$ws.Range("M76").Value = -10308.625
$ws.Range("N76").Value = -7630
$ws.Range("H79").Value = 10221
$ws.Range("I79").Value = 10623.625
$ws.Range("J79").Value = 7000
$ws.Range("K79").Value = 10623.625
$ws.Range("L79").Value = 7000
$ws.Range("M79").Value = -9531.625
$ws.Range("N79").Value = -9184
$ws.Range("H86").Value = 1669.8667
$ws.Range("I86").Value = 2338.7144
$ws.Range("J86").Value = 1084.625
$ws.Range("K86").Value = 2338.7144
$ws.Range("L86").Value = 1084.625
$ws.Range("M86").Value = -1215.7144
$ws.Range("N86").Value = -3330.625
$ws.Range("H89").Value = 1669.8667
$ws.Range("I89").Value = 2338.7144
$ws.Range("J89").Value = 1084.625
$ws.Range("K89").Value = 11693.572
$ws.Range("L89").Value = 5423.125
$ws.Range("M89").Value = -6077.572
$ws.Range("N89").Value = -16655.125
$ws.Range("H92").Value = 43445.25
$ws.Range("I92").Value = 72713.14
$ws.Range("J92").Value = 2470.2
$ws.Range("K92").Value = 72713.14
$ws.Range("L92").Value = 2470.2
$ws.Range("M92").Value = -71465.14
$ws.Range("N92").Value = -4966.2
$ws.Range("H97").Value = 952.5
$ws.Range("J97").Value = 952.5
$ws.Range("L97").Value = 2857.5
$ws.Range("N97").Value = -3849.5
$ws.Range("H99").Value = 330.81818
$ws.Range("I99").Value = 340.55554
$ws.Range("J99").Value = 287
$ws.Range("K99").Value = 1021.66662
$ws.Range("L99").Value = 861
$ws.Range("M99").Value = 476.33338
$ws.Range("N99").Value = -3857
$ws.Range("H101").Value = 1760.7273
$ws.Range("I101").Value = 1143
$ws.Range("J101").Value = 3408
$ws.Range("K101").Value = 3429
$ws.Range("L101").Value = 10224
$ws.Range("M101").Value = -1807
$ws.Range("N101").Value = -13468
$ws.Range("H107").Value = 1736
$ws.Range("I107").Value = 1157.4375
$ws.Range("K107").Value = 1157.4375
$ws.Range("M107").Value = 762.5625
$ws.Range("H111").Value = 2499
$ws.Range("I111").Value = 2499
$ws.Range("J111").Value = 2499
$ws.Range("K111").Value = 7497
$ws.Range("L111").Value = 7497
$ws.Range("M111").Value = -4430
$ws.Range("N111").Value = -13631
$ws.Range("H113").Value = 5329.316
$ws.Range("I113").Value = 6457.6665
$ws.Range("J113").Value = 3395
$ws.Range("K113").Value = 6457.6665
$ws.Range("L113").Value = 3395
$ws.Range("M113").Value = -3203.6665
$ws.Range("N113").Value = -9903
$ws.Range("H116").Value = 7023.0293
$ws.Range("J116").Value = 8124.25
$ws.Range("L116").Value = 8124.25
$ws.Range("N116").Value = -15008.25
$ws.Range("H125").Value = 3926.6667
$ws.Range("I125").Value = 3923
$ws.Range("K125").Value = 35307
$ws.Range("M125").Value = -32847
$ws.Range("H132").Value = 2436.6765
$ws.Range("I132").Value = 2190.4614
$ws.Range("K132").Value = 6571.3842
$ws.Range("M132").Value = -4041.3842
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = $null
$ws.Range("H137").Value = 354756.25
$ws.Range("I137").Value = 2094.6924
$ws.Range("K137").Value = 6284.0772
$ws.Range("M137").Value = -3734.0772
$ws.Range("H138").Value = 2530.32
$ws.Range("I138").Value = 1295.2258
$ws.Range("J138").Value = 3400.5
$ws.Range("K138").Value = 3885.6774
$ws.Range("L138").Value = 10201.5
$ws.Range("M138").Value = 1254.3226
$ws.Range("N138").Value = -20481.5
$ws.Range("H139").Value = 73146.336
$ws.Range("J139").Value = 73146.336
$ws.Range("L139").Value = 73146.336
$ws.Range("N139").Value = -83426.336
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = $null
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1923.5161
$ws.Range("I2").Value = 1961.8214
$ws.Range("K2").Value = 1961.8214
$ws.Range("M2").Value = -1848.8214
$ws.Range("H32").Value = 5188.2656
$ws.Range("I32").Value = 4001.4834
$ws.Range("K32").Value = 4001.4834
$ws.Range("M32").Value = -3714.4834
$ws.Range("H45").Value = 2021.7587
$ws.Range("I45").Value = 1031
$ws.Range("J45").Value = 6777.4
$ws.Range("K45").Value = 1031
$ws.Range("L45").Value = 6777.4
$ws.Range("M45").Value = -654
$ws.Range("N45").Value = -7531.4
$ws.Range("H74").Value = 60195.324
$ws.Range("I74").Value = 36571.77
$ws.Range("K74").Value = 36571.77
$ws.Range("M74").Value = -35697.77
$ws.Range("H77").Value = 60195.324
$ws.Range("I77").Value = 36571.77
$ws.Range("K77").Value = 182858.85
$ws.Range("M77").Value = -178490.85
$ws.Range("H102").Value = 5937.421
$ws.Range("I102").Value = 5808.769
$ws.Range("J102").Value = 6216.1665
$ws.Range("K102").Value = 5808.769
$ws.Range("L102").Value = 6216.1665
$ws.Range("M102").Value = -4186.769
$ws.Range("N102").Value = -9460.166499999999
$ws.Range("H116").Value = 1923.5161
$ws.Range("I116").Value = 1961.8214
$ws.Range("K116").Value = 1961.8214
$ws.Range("M116").Value = 332.1786
$ws.Range("H122").Value = 48415.7
$ws.Range("I122").Value = 2539.3845
$ws.Range("J122").Value = 133614.58
$ws.Range("K122").Value = 7618.1535
$ws.Range("L122").Value = 400843.74
$ws.Range("M122").Value = -5168.1535
$ws.Range("N122").Value = -405743.74
$ws.Range("H132").Value = 2307.2
$ws.Range("I132").Value = 2221.3
$ws.Range("J132").Value = 2822.6
$ws.Range("K132").Value = 6663.900000000001
$ws.Range("L132").Value = 8467.799999999999
$ws.Range("M132").Value = -4133.900000000001
$ws.Range("N132").Value = -13527.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1923.5161
$ws.Range("I3").Value = 1961.8214
$ws.Range("K3").Value = 1961.8214
$ws.Range("M3").Value = -1847.8214
$ws.Range("H86").Value = 57156.062
$ws.Range("I86").Value = 945.63635
$ws.Range("J86").Value = 180819
$ws.Range("K86").Value = 945.63635
$ws.Range("L86").Value = 180819
$ws.Range("M86").Value = 177.36365
$ws.Range("N86").Value = -183065
$ws.Range("H89").Value = 57156.062
$ws.Range("I89").Value = 945.63635
$ws.Range("J89").Value = 180819
$ws.Range("K89").Value = 4728.18175
$ws.Range("L89").Value = 904095
$ws.Range("M89").Value = 887.8182500000003
$ws.Range("N89").Value = -915327
$ws.Range("H94").Value = 60861.8
$ws.Range("I94").Value = 770
$ws.Range("J94").Value = 150999.5
$ws.Range("K94").Value = 770
$ws.Range("L94").Value = 150999.5
$ws.Range("M94").Value = -319
$ws.Range("N94").Value = -151901.5
$ws.Range("H99").Value = 1810.826
$ws.Range("I99").Value = 1542
$ws.Range("K99").Value = 1542
$ws.Range("M99").Value = -44
$ws.Range("H105").Value = 2086.7222
$ws.Range("I105").Value = 1851.4
$ws.Range("J105").Value = 3263.3333
$ws.Range("K105").Value = 1851.4
$ws.Range("L105").Value = 3263.3333
$ws.Range("M105").Value = -104.4000000000001
$ws.Range("N105").Value = -6757.3333
$ws.Range("H107").Value = 2240.9
$ws.Range("I107").Value = 2484.25
$ws.Range("J107").Value = 1267.5
$ws.Range("K107").Value = 2484.25
$ws.Range("L107").Value = 1267.5
$ws.Range("M107").Value = -564.25
$ws.Range("N107").Value = -5107.5
$ws.Range("H134").Value = 4108.222
$ws.Range("I134").Value = 2409.12
$ws.Range("J134").Value = 7969.8184
$ws.Range("K134").Value = 7227.36
$ws.Range("L134").Value = 23909.4552
$ws.Range("M134").Value = -4692.36
$ws.Range("N134").Value = -28979.4552
$ws.Range("H138").Value = 96658.414
$ws.Range("J138").Value = 96658.414
$ws.Range("L138").Value = 96658.414
$ws.Range("N138").Value = -106938.414

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 224.84616
$ws.Range("I7").Value = 109.25
$ws.Range("J7").Value = 409.8
$ws.Range("K7").Value = 109.25
$ws.Range("L7").Value = 409.8
$ws.Range("M7").Value = 3.75
$ws.Range("N7").Value = -635.8
$ws.Range("H31").Value = 2270.7666
$ws.Range("I31").Value = 1565.8
$ws.Range("J31").Value = 2975.7334
$ws.Range("K31").Value = 1565.8
$ws.Range("L31").Value = 2975.7334
$ws.Range("M31").Value = -1270.8
$ws.Range("N31").Value = -3565.7334
$ws.Range("H34").Value = 2270.7666
$ws.Range("I34").Value = 1565.8
$ws.Range("J34").Value = 2975.7334
$ws.Range("K34").Value = 1565.8
$ws.Range("L34").Value = 2975.7334
$ws.Range("M34").Value = -1363.8
$ws.Range("N34").Value = -3379.7334
$ws.Range("H97").Value = 49374.5
$ws.Range("J97").Value = 49998.668
$ws.Range("L97").Value = 49998.668
$ws.Range("N97").Value = -51980.668
$ws.Range("H99").Value = 5345.304
$ws.Range("I99").Value = 5160.6665
$ws.Range("J99").Value = 5691.5
$ws.Range("K99").Value = 5160.6665
$ws.Range("L99").Value = 5691.5
$ws.Range("M99").Value = -3662.6665
$ws.Range("N99").Value = -8687.5
$ws.Range("H107").Value = 28081.229
$ws.Range("I107").Value = 42353.773
$ws.Range("J107").Value = 3927.6924
$ws.Range("K107").Value = 42353.773
$ws.Range("L107").Value = 3927.6924
$ws.Range("M107").Value = -40433.773
$ws.Range("N107").Value = -7767.6924
$ws.Range("H122").Value = 1525.8889
$ws.Range("J122").Value = 1679.5
$ws.Range("L122").Value = 5038.5
$ws.Range("N122").Value = -9938.5
$ws.Range("H126").Value = 5345.304
$ws.Range("I126").Value = 5160.6665
$ws.Range("J126").Value = 5691.5
$ws.Range("K126").Value = 15481.9995
$ws.Range("L126").Value = 17074.5
$ws.Range("M126").Value = -13011.9995
$ws.Range("N126").Value = -22014.5
$ws.Range("H132").Value = 7296.1665
$ws.Range("I132").Value = 2143.8333
$ws.Range("K132").Value = 6431.499899999999
$ws.Range("M132").Value = -3901.499899999999
$ws.Range("H134").Value = 1310.697
$ws.Range("I134").Value = 1305.275
$ws.Range("J134").Value = 1319.0385
$ws.Range("K134").Value = 3915.825
$ws.Range("L134").Value = 3957.1155
$ws.Range("M134").Value = -1380.825
$ws.Range("N134").Value = -9027.1155
$ws.Range("H141").Value = 175013.03
$ws.Range("J141").Value = 181453.56
$ws.Range("L141").Value = 181453.56
$ws.Range("N141").Value = -191813.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 409.7143
$ws.Range("I23").Value = 80.75
$ws.Range("J23").Value = 848.3333
$ws.Range("K23").Value = 242.25
$ws.Range("L23").Value = 2544.9999
$ws.Range("M23").Value = -7.25
$ws.Range("N23").Value = -3014.9999
$ws.Range("H40").Value = 124.23077
$ws.Range("I40").Value = 106.9
$ws.Range("J40").Value = 182
$ws.Range("K40").Value = 427.6
$ws.Range("L40").Value = 728
$ws.Range("M40").Value = -358.6
$ws.Range("N40").Value = -866
$ws.Range("H43").Value = 7499
$ws.Range("I43").Value = 7998
$ws.Range("J43").Value = 7000
$ws.Range("K43").Value = 23994
$ws.Range("L43").Value = 21000
$ws.Range("M43").Value = -23880
$ws.Range("N43").Value = -21228
$ws.Range("H81").Value = 4499.5
$ws.Range("J81").Value = 4499.5
$ws.Range("L81").Value = 13498.5
$ws.Range("N81").Value = -15744.5
$ws.Range("H84").Value = 4499.5
$ws.Range("J84").Value = 4499.5
$ws.Range("L84").Value = 40495.5
$ws.Range("N84").Value = -51727.5
$ws.Range("H114").Value = 3029.7778
$ws.Range("I114").Value = 2609.7144
$ws.Range("K114").Value = 7829.1432
$ws.Range("M114").Value = -4575.1432
$ws.Range("H129").Value = 4302.08
$ws.Range("J129").Value = 5156.2
$ws.Range("L129").Value = 15468.6
$ws.Range("N129").Value = -25468.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 45107
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null
$ws.Range("H80").Value = 14701.333
$ws.Range("I80").Value = 7249
$ws.Range("J80").Value = 16830.572
$ws.Range("K80").Value = 7249
$ws.Range("L80").Value = 16830.572
$ws.Range("M80").Value = -6251
$ws.Range("N80").Value = -18826.572
$ws.Range("H81").Value = 45107
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null
$ws.Range("H83").Value = 14701.333
$ws.Range("I83").Value = 7249
$ws.Range("J83").Value = 16830.572
$ws.Range("K83").Value = 36245
$ws.Range("L83").Value = 84152.86
$ws.Range("M83").Value = -31253
$ws.Range("N83").Value = -94136.86
$ws.Range("H84").Value = 45107
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null
$ws.Range("H94").Value = 47218.168
$ws.Range("J94").Value = 47218.168
$ws.Range("L94").Value = 47218.168
$ws.Range("N94").Value = -48570.168
$ws.Range("H97").Value = 4001.5557
$ws.Range("I97").Value = 4000
$ws.Range("J97").Value = 4004.6667
$ws.Range("K97").Value = 4000
$ws.Range("L97").Value = 4004.6667
$ws.Range("M97").Value = -3504
$ws.Range("N97").Value = -4996.6667
$ws.Range("H102").Value = 10224.412
$ws.Range("I102").Value = 8419.392
$ws.Range("J102").Value = 13998.546
$ws.Range("K102").Value = 8419.392
$ws.Range("L102").Value = 13998.546
$ws.Range("M102").Value = -6797.392
$ws.Range("N102").Value = -17242.546
$ws.Range("H107").Value = 184.66667
$ws.Range("J107").Value = 195.5
$ws.Range("L107").Value = 195.5
$ws.Range("N107").Value = -4035.5
$ws.Range("H109").Value = 31248.75
$ws.Range("J109").Value = 31248.75
$ws.Range("L109").Value = 31248.75
$ws.Range("N109").Value = -33328.75
$ws.Range("H113").Value = 14561.308
$ws.Range("J113").Value = 6449.5
$ws.Range("L113").Value = 6449.5
$ws.Range("N113").Value = -10789.5
$ws.Range("H122").Value = 3051.077
$ws.Range("I122").Value = 1296.5555
$ws.Range("K122").Value = 3889.6665
$ws.Range("M122").Value = -1439.6665
$ws.Range("H132").Value = 3505.6
$ws.Range("I132").Value = 3633.2927
$ws.Range("J132").Value = 2196.75
$ws.Range("K132").Value = 10899.8781
$ws.Range("L132").Value = 6590.25
$ws.Range("M132").Value = -8369.8781
$ws.Range("N132").Value = -11650.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 31585.5
$ws.Range("I7").Value = 36182.668
$ws.Range("J7").Value = 4002.5
$ws.Range("K7").Value = 36182.668
$ws.Range("L7").Value = 4002.5
$ws.Range("M7").Value = -36070.668
$ws.Range("N7").Value = -4226.5
$ws.Range("H22").Value = 7288.1562
$ws.Range("I22").Value = 3887.25
$ws.Range("J22").Value = 7774
$ws.Range("K22").Value = 3887.25
$ws.Range("L22").Value = 7774
$ws.Range("M22").Value = -3592.25
$ws.Range("N22").Value = -8364
$ws.Range("H27").Value = 7288.1562
$ws.Range("I27").Value = 3887.25
$ws.Range("J27").Value = 7774
$ws.Range("K27").Value = 3887.25
$ws.Range("L27").Value = 7774
$ws.Range("M27").Value = -3780.25
$ws.Range("N27").Value = -7988
$ws.Range("H40").Value = 27976.143
$ws.Range("I40").Value = 30000.562
$ws.Range("K40").Value = 30000.562
$ws.Range("M40").Value = -29864.562
$ws.Range("H46").Value = 1927.7059
$ws.Range("J46").Value = 1942.5625
$ws.Range("L46").Value = 1942.5625
$ws.Range("N46").Value = -2318.5625
$ws.Range("H61").Value = 4075.818
$ws.Range("I61").Value = 4003.7778
$ws.Range("J61").Value = 4400
$ws.Range("K61").Value = 4003.7778
$ws.Range("L61").Value = 4400
$ws.Range("M61").Value = -3801.7778
$ws.Range("N61").Value = -4804
$ws.Range("H68").Value = 3676.1177
$ws.Range("I68").Value = 2076.4614
$ws.Range("K68").Value = 2076.4614
$ws.Range("M68").Value = -1327.4614
$ws.Range("H71").Value = 3676.1177
$ws.Range("I71").Value = 2076.4614
$ws.Range("K71").Value = 10382.307
$ws.Range("M71").Value = -6638.307000000001
$ws.Range("H93").Value = 1286.625
$ws.Range("I93").Value = 519.6
$ws.Range("J93").Value = 2565
$ws.Range("K93").Value = 519.6
$ws.Range("L93").Value = 2565
$ws.Range("M93").Value = 728.4
$ws.Range("N93").Value = -5061
$ws.Range("H94").Value = 69998.5
$ws.Range("J94").Value = 69998.5
$ws.Range("L94").Value = 69998.5
$ws.Range("N94").Value = -71350.5
$ws.Range("H100").Value = 2368.182
$ws.Range("I100").Value = 2564.8333
$ws.Range("J100").Value = 2132.2
$ws.Range("K100").Value = 2564.8333
$ws.Range("L100").Value = 2132.2
$ws.Range("M100").Value = -2023.8333
$ws.Range("N100").Value = -3214.2
$ws.Range("H113").Value = 4075.818
$ws.Range("I113").Value = 4003.7778
$ws.Range("J113").Value = 4400
$ws.Range("K113").Value = 4003.7778
$ws.Range("L113").Value = 4400
$ws.Range("M113").Value = -1833.7778
$ws.Range("N113").Value = -8740
$ws.Range("H122").Value = 5817.2188
$ws.Range("I122").Value = 3362.6072
$ws.Range("J122").Value = 22999.5
$ws.Range("K122").Value = 10087.8216
$ws.Range("L122").Value = 68998.5
$ws.Range("M122").Value = -7637.821599999999
$ws.Range("N122").Value = -73898.5
$ws.Range("H126").Value = 31585.5
$ws.Range("I126").Value = 36182.668
$ws.Range("J126").Value = 4002.5
$ws.Range("K126").Value = 108548.004
$ws.Range("L126").Value = 12007.5
$ws.Range("M126").Value = -106078.004
$ws.Range("N126").Value = -16947.5
$ws.Range("H132").Value = 4486.3887
$ws.Range("I132").Value = 2946.4167
$ws.Range("K132").Value = 8839.250100000001
$ws.Range("M132").Value = -6309.250100000001
$ws.Range("H136").Value = 4482.15
$ws.Range("I136").Value = 2903.75
$ws.Range("K136").Value = 8711.25
$ws.Range("M136").Value = -6161.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 69995
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
$ws.Range("H67").Value = 69995
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
$ws.Range("H96").Value = 7049.5
$ws.Range("I96").Value = 8299.4
$ws.Range("K96").Value = 8299.4
$ws.Range("M96").Value = -6926.4
$ws.Range("H107").Value = 2100
$ws.Range("I107").Value = 1416.6666
$ws.Range("K107").Value = 4249.9998
$ws.Range("M107").Value = -2329.9998
$ws.Range("H122").Value = 3603.8462
$ws.Range("I122").Value = 3897.7058
$ws.Range("J122").Value = 3048.7778
$ws.Range("K122").Value = 11693.1174
$ws.Range("L122").Value = 9146.3334
$ws.Range("M122").Value = -9243.117400000001
$ws.Range("N122").Value = -14046.3334
$ws.Range("H132").Value = 3521.9412
$ws.Range("I132").Value = 1867.0938
$ws.Range("K132").Value = 5601.2814
$ws.Range("M132").Value = -3071.2814
$ws.Range("H136").Value = 3568.1282
$ws.Range("I136").Value = 3091.0625
$ws.Range("K136").Value = 9273.1875
$ws.Range("M136").Value = -6723.1875
